$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gumottc_characteristics")

# --- Row 26 becomes the new "elevation_estimate" row. It is no longer a
# group header (that role moves to row 27 now), so it loses the bottom
# border that previously marked the start of a new field-group and instead
# takes on the plain "interior row" look (copy format from a known plain
# cell, e.g. C1).
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A26:E26").PasteSpecial(-4122) | Out-Null

$ws.Range("A26").Value = "Y"
$ws.Range("B26").Value = "elevation_estimate"
$ws.Range("C26").Value = "if elevation had to be estimated using coordinates determined via WGS84"
$ws.Range("D26").Value = "n"
$ws.Range("E26").Value = "NA"

# --- Row 27 becomes "match_below_500ft" and is now the group header, so
# column A picks up the bottom-border "new group" styling (copy from an
# existing header cell, e.g. A24). The rest of the row keeps the plain look.
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A27").PasteSpecial(-4122) | Out-Null

$ws.Range("A27").Value = "Z"
$ws.Range("B27").Value = "match_below_500ft"
$ws.Range("C27").Value = "if sample has a species match in the other time period that was collected within 500ft of each other in terms of elevation"
$ws.Range("D27").Value = "n"
$ws.Range("E27").Value = "NA"

# --- Row 28 is a brand new row for "match_500ft_id" (pushed down from the
# old row 27). Column A is the new group header for this field; the rest
# of the row is plain.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A28:E28").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null

$ws.Range("A28").Value = "AA"
$ws.Range("B28").Value = "match_500ft_id"
$ws.Range("C28").Value = "ID code of sample collected within 500 ft of its matched sample if applicable"
$ws.Range("D28").Value = "n"
$ws.Range("E28").Value = "NA"
